# Demo Verification Script and Display CF
# Update the "Date" column (column B) timestamps on the three worksheets
# to reflect the latest test run.

$wb = $excel.ActiveWorkbook

# --- CC-Payments-Sale ---
$ws1 = $wb.Worksheets.Item("CC-Payments-Sale")
$ws1.Range("B2").Value = "Thu May 01 19:12:53 IST 2025"
$ws1.Range("B3").Value = "Thu May 01 19:13:35 IST 2025"
$ws1.Range("B4").Value = "Thu May 01 19:14:11 IST 2025"
$ws1.Range("B5").Value = "Thu May 01 19:14:44 IST 2025"
$ws1.Range("B6").Value = "Thu May 01 19:15:21 IST 2025"
$ws1.Range("B7").Value = "Thu May 01 19:15:54 IST 2025"

# --- CC-Payments-Auth ---
$ws2 = $wb.Worksheets.Item("CC-Payments-Auth")
$ws2.Range("B2").Value = "Thu May 01 19:05:59 IST 2025"
$ws2.Range("B3").Value = "Thu May 01 19:06:45 IST 2025"
$ws2.Range("B4").Value = "Thu May 01 19:07:23 IST 2025"
$ws2.Range("B5").Value = "Thu May 01 19:07:59 IST 2025"
$ws2.Range("B6").Value = "Thu May 01 19:08:37 IST 2025"
$ws2.Range("B7").Value = "Thu May 01 19:09:13 IST 2025"

# --- ACH-Payments-Debit ---
$ws3 = $wb.Worksheets.Item("ACH-Payments-Debit")
$ws3.Range("B8").Value = "Thu May 01 19:09:50 IST 2025"
$ws3.Range("B9").Value = "Thu May 01 19:10:48 IST 2025"
$ws3.Range("B10").Value = "Thu May 01 19:11:59 IST 2025"
